$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Actualizacion desde MV -datos-": append the latest daily rate rows
# (Serie / 2 anios / 5 anios / 10 anios) after the last existing row (175).
$newRows = @(
    @{ Fecha = "13-09-2021"; C2 = $null; C5 = 3.43; C10 = 3.26 },
    @{ Fecha = "14-09-2021"; C2 = 4.51;  C5 = 3.54; C10 = 3.25 },
    @{ Fecha = "15-09-2021"; C2 = 4.53;  C5 = 3.71; C10 = 3.27 },
    @{ Fecha = "16-09-2021"; C2 = 4.62;  C5 = $null; C10 = 3.26 }
)

$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($row, 1).Value = $data.Fecha
    if ($null -ne $data.C2) {
        $ws.Cells.Item($row, 3).Value = $data.C2
    }
    if ($null -ne $data.C5) {
        $ws.Cells.Item($row, 4).Value = $data.C5
    }
    if ($null -ne $data.C10) {
        $ws.Cells.Item($row, 5).Value = $data.C10
    }
}
